$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 183.42857
$ws.Range("I33").Value = 106.42857
$ws.Range("J33").Value = 337.42856
$ws.Range("K33").Value = 106.42857
$ws.Range("L33").Value = 337.42856
$ws.Range("M33").Value = 122.57143
$ws.Range("N33").Value = -795.4285600000001
$ws.Range("H41").Value = 599.2
$ws.Range("I41").Value = 248
$ws.Range("J41").Value = 687
$ws.Range("K41").Value = 248
$ws.Range("L41").Value = 687
$ws.Range("M41").Value = 192
$ws.Range("N41").Value = -1567
$ws.Range("H107").Value = 882.8946999999999
$ws.Range("I107").Value = 882.8946999999999
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 882.8946999999999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1037.1053
$ws.Range("N107").ClearContents()
$ws.Range("H108").Value = 32329.5
$ws.Range("J108").Value = 32329.5
$ws.Range("L108").Value = 32329.5
$ws.Range("N108").Value = -40009.5
$ws.Range("H109").Value = 38680
$ws.Range("J109").Value = 38680
$ws.Range("L109").Value = 38680
$ws.Range("N109").Value = -41454
$ws.Range("H128").Value = 51240.668
$ws.Range("J128").Value = 51240.668
$ws.Range("L128").Value = 51240.668
$ws.Range("N128").Value = -61200.668
$ws.Range("H133").Value = 73884.86
$ws.Range("J133").Value = 73884.86
$ws.Range("L133").Value = 73884.86
$ws.Range("N133").Value = -84004.86
$ws.Range("H136").Value = 55049.8
$ws.Range("J136").Value = 55049.8
$ws.Range("L136").Value = 55049.8
$ws.Range("N136").Value = -65249.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 23500
$ws.Range("I34").Value = 15000
$ws.Range("J34").Value = 32000
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 32000
$ws.Range("M34").Value = -14729
$ws.Range("N34").Value = -32542
$ws.Range("H61").Value = 1991.1136
$ws.Range("I61").Value = 1860.25
$ws.Range("J61").Value = 3299.75
$ws.Range("K61").Value = 1860.25
$ws.Range("L61").Value = 3299.75
$ws.Range("M61").Value = -1648.25
$ws.Range("N61").Value = -3723.75
$ws.Range("H107").Value = 38998.332
$ws.Range("J107").Value = 38998.332
$ws.Range("L107").Value = 38998.332
$ws.Range("N107").Value = -46678.332
$ws.Range("H109").Value = 41251.332
$ws.Range("J109").Value = 41251.332
$ws.Range("L109").Value = 41251.332
$ws.Range("N109").Value = -44025.332
$ws.Range("H117").Value = 44229.4
$ws.Range("J117").Value = 44229.4
$ws.Range("L117").Value = 44229.4
$ws.Range("N117").Value = -53407.4
$ws.Range("H118").Value = 49650.5
$ws.Range("J118").Value = 49650.5
$ws.Range("L118").Value = 49650.5
$ws.Range("N118").Value = -52964.5
$ws.Range("H123").Value = 38214.5
$ws.Range("J123").Value = 51429
$ws.Range("L123").Value = 51429
$ws.Range("N123").Value = -61229
$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("N125").Value = -60555
$ws.Range("H130").Value = 43215.8
$ws.Range("J130").Value = 43215.8
$ws.Range("L130").Value = 43215.8
$ws.Range("N130").Value = -53255.8
$ws.Range("H131").Value = 51235.668
$ws.Range("J131").Value = 51235.668
$ws.Range("L131").Value = 51235.668
$ws.Range("N131").Value = -61315.668
$ws.Range("H136").Value = 1991.1136
$ws.Range("I136").Value = 1860.25
$ws.Range("J136").Value = 3299.75
$ws.Range("K136").Value = 5580.75
$ws.Range("L136").Value = 9899.25
$ws.Range("M136").Value = -3030.75
$ws.Range("N136").Value = -14999.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 49744.668
$ws.Range("J117").Value = 49744.668
$ws.Range("L117").Value = 49744.668
$ws.Range("N117").Value = -58922.668
$ws.Range("H119").Value = 48761
$ws.Range("J119").Value = 48761
$ws.Range("L119").Value = 48761
$ws.Range("N119").Value = -58437
$ws.Range("H124").Value = 49996
$ws.Range("J124").Value = 49996
$ws.Range("L124").Value = 49996
$ws.Range("N124").Value = -59816
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H126").Value = 50780
$ws.Range("J126").Value = 50780
$ws.Range("L126").Value = 50780
$ws.Range("N126").Value = -60660
$ws.Range("H130").Value = 47306.668
$ws.Range("J130").Value = 47306.668
$ws.Range("L130").Value = 47306.668
$ws.Range("N130").Value = -57346.668
$ws.Range("H134").Value = 2351.375
$ws.Range("I134").Value = 1958.9231
$ws.Range("J134").Value = 4052
$ws.Range("K134").Value = 5876.7693
$ws.Range("L134").Value = 12156
$ws.Range("M134").Value = -3341.7693
$ws.Range("N134").Value = -17226
$ws.Range("H137").Value = 56741.332
$ws.Range("J137").Value = 56741.332
$ws.Range("L137").Value = 56741.332
$ws.Range("N137").Value = -66941.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 40283.168
$ws.Range("J20").Value = 40283.168
$ws.Range("L20").Value = 40283.168
$ws.Range("N20").Value = -40755.168
$ws.Range("H30").Value = 40283.168
$ws.Range("J30").Value = 40283.168
$ws.Range("L30").Value = 40283.168
$ws.Range("N30").Value = -40465.168
$ws.Range("H116").Value = 44350
$ws.Range("J116").Value = 44350
$ws.Range("L116").Value = 44350
$ws.Range("N116").Value = -53528
$ws.Range("H128").Value = 40283.168
$ws.Range("J128").Value = 40283.168
$ws.Range("L128").Value = 40283.168
$ws.Range("N128").Value = -50243.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 37040916
$ws.Range("I141").Value = 52634990
$ws.Range("J141").Value = 4993.5
$ws.Range("K141").Value = 157904970
$ws.Range("L141").Value = 14980.5
$ws.Range("M141").Value = -157899790
$ws.Range("N141").Value = -25340.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 42420
$ws.Range("J121").Value = 42420
$ws.Range("L121").Value = 42420
$ws.Range("N121").Value = -45914
$ws.Range("H127").Value = 47909.57
$ws.Range("J127").Value = 47909.57
$ws.Range("L127").Value = 47909.57
$ws.Range("N127").Value = -57829.57
$ws.Range("H130").Value = 47441.668
$ws.Range("J130").Value = 47441.668
$ws.Range("L130").Value = 47441.668
$ws.Range("N130").Value = -57481.668
$ws.Range("H136").Value = 1613.7906
$ws.Range("I136").Value = 1318.7949
$ws.Range("K136").Value = 3956.384700000001
$ws.Range("M136").Value = -1406.384700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 42480.8
$ws.Range("J120").Value = 42480.8
$ws.Range("L120").Value = 42480.8
$ws.Range("N120").Value = -52156.8
$ws.Range("H135").Value = 57128.285
$ws.Range("J135").Value = 57128.285
$ws.Range("L135").Value = 57128.285
$ws.Range("N135").Value = -67268.285
$ws.Range("H136").Value = 17257.662
$ws.Range("I136").Value = 40825.92
$ws.Range("J136").Value = 2527.5
$ws.Range("K136").Value = 122477.76
$ws.Range("L136").Value = 7582.5
$ws.Range("M136").Value = -119927.76
$ws.Range("N136").Value = -12682.5
